$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 53

$ws.Cells.Item($row, 1).Value = 8
$ws.Cells.Item($row, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item($row, 3).Value = "Coquimbo"

$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($row, 4).Value = (Get-Date -Year 2021 -Month 11 -Day 9 -Hour 0 -Minute 0 -Second 0)

$ws.Cells.Item($row, 5).Value = 4
$ws.Cells.Item($row, 6).Value = 100112052
$ws.Cells.Item($row, 7).Value = "Albahaca"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 800
$ws.Cells.Item($row, 11).Value = 3500
$ws.Cells.Item($row, 12).Value = 4000
$ws.Cells.Item($row, 13).Value = 3750
$ws.Cells.Item($row, 14).Value = "$/paquete"
$ws.Cells.Item($row, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item($row, 16).Value = 3750
$ws.Cells.Item($row, 17).Value = 1
$ws.Cells.Item($row, 18).Value = "Hortaliza"
